# Applies the "Updated cryptos list" price/volume refresh to Sheet1.
# Numeric-looking Price-column values are written with a leading "'"
# (Excel quote-prefix) so they stay plain text, matching the source data
# (e.g. "1.00", "0.0000189") instead of being auto-parsed into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "73.846.93"
$ws.Range("E2").Value = "  +7.37%  "
# Row 3
$ws.Range("D3").Value = "2.617.93"
$ws.Range("E3").Value = "  +7.18%  "
# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.04%  "
# Row 5
$ws.Range("D5").Value = "'186.64"
$ws.Range("E5").Value = "  +14.10%  "
# Row 6
$ws.Range("D6").Value = "'580.82"
$ws.Range("E6").Value = "  +3.57%  "
# Row 7
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.06%  "
# Row 8
$ws.Range("E8").Value = "  +4.55%  "
# Row 9
$ws.Range("D9").Value = "'0.198"
$ws.Range("E9").Value = "  +16.18%  "
# Row 10
$ws.Range("D10").Value = "2.615.98"
$ws.Range("E10").Value = "  +7.11%  "
# Row 11
$ws.Range("E11").Value = "  +1.16%  "
# Row 12
$ws.Range("D12").Value = "'0.357"
$ws.Range("E12").Value = "  +7.48%  "
# Row 13
$ws.Range("E13").Value = "  +1.42%  "
# Row 14
$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").Value = "74.135.83"
$ws.Range("E14").Value = "  +7.99%  "
# Row 15
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000189"
$ws.Range("E15").Value = "  +5.24%  "
# Row 16
$ws.Range("E16").Value = "  +7.34%  "
# Row 17
$ws.Range("D17").Value = "'26.37"
$ws.Range("E17").Value = "  +12.61%  "
# Row 18
$ws.Range("D18").Value = "2.611.47"
$ws.Range("E18").Value = "  +6.93%  "
# Row 19
$ws.Range("E19").Value = "  +29.36%  "
# Row 20
$ws.Range("D20").Value = "'11.80"
$ws.Range("E20").Value = "  +11.34%  "
# Row 21
$ws.Range("D21").Value = "'366.27"
$ws.Range("E21").Value = "  +7.94%  "
# Row 22
$ws.Range("D22").Value = "'2.29"
$ws.Range("E22").Value = "  +17.64%  "
# Row 23
$ws.Range("E23").Value = "  +5.76%  "
# Row 24
$ws.Range("E24").Value = "  -0.10%  "
# Row 25
$ws.Range("D25").Value = "'69.83"
$ws.Range("E25").Value = "  +6.60%  "
# Row 26
$ws.Range("D26").Value = "'4.12"
$ws.Range("E26").Value = "  +8.68%  "
# Row 27
$ws.Range("D27").Value = "'9.30"
$ws.Range("E27").Value = "  +10.91%  "
# Row 28
$ws.Range("D28").Value = "2.756.63"
$ws.Range("E28").Value = "  +7.34%  "
# Row 29
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.14%  "
# Row 30
$ws.Range("D30").Value = "0.0₃0940"
$ws.Range("E30").Value = "  +13.87%  "
# Row 31
$ws.Range("D31").Value = "'520.64"
$ws.Range("E31").Value = "  +20.20%  "
# Row 32
$ws.Range("E32").Value = "  +14.79%  "
# Row 33
$ws.Range("D33").Value = "'7.64"
$ws.Range("E33").Value = "  +6.23%  "
# Row 34
$ws.Range("D34").Value = "'1.74"
$ws.Range("E34").Value = "  +8.82%  "
# Row 35
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.01%  "
# Row 36
$ws.Range("D36").Value = "'162.61"
$ws.Range("E36").Value = "  +2.06%  "
# Row 37
$ws.Range("E37").Value = "  +10.43%  "
# Row 38
$ws.Range("D38").Value = "'19.10"
$ws.Range("E38").Value = "  +6.08%  "
# Row 39
$ws.Range("D39").Value = "'19.26"
$ws.Range("E39").Value = "  +1.40%  "
# Row 40
$ws.Range("E40").Value = "  +0.08%  "
# Row 41
$ws.Range("E41").Value = "  +11.82%  "
# Row 42
$ws.Range("E42").Value = "  +9.26%  "
# Row 43
$ws.Range("E43").Value = "  +7.90%  "
# Row 44
$ws.Range("D44").Value = "'161.32"
$ws.Range("E44").Value = "  +24.13%  "
# Row 45
$ws.Range("D45").Value = "'2.37"
$ws.Range("E45").Value = "  +13.86%  "
# Row 46
$ws.Range("D46").Value = "'1.17"
$ws.Range("E46").Value = "  +8.79%  "
# Row 47
$ws.Range("D47").Value = "'38.89"
$ws.Range("E47").Value = "  +3.58%  "
# Row 48
$ws.Range("D48").Value = "'0.0849"
$ws.Range("E48").Value = "  +18.01%  "
# Row 49
$ws.Range("D49").Value = "'3.60"
$ws.Range("E49").Value = "  +8.01%  "
# Row 50
$ws.Range("D50").Value = "'0.523"
$ws.Range("E50").Value = "  +7.91%  "
# Row 51
$ws.Range("D51").Value = "'20.69"
$ws.Range("E51").Value = "  +21.80%  "
